$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows for line7/line8 right after the line6 row (row 7),
# before the extr1..extr8 rows (which were rows 8-15, now become 10-17).
$ws.Rows("8:9").Insert()

# Copy formatting (bold border style) from the row above (row 7's column A)
# into the new A8:A9 cells so they match the existing styled "index" column.
$ws.Range("A7").Copy()
$ws.Range("A8:A9").PasteSpecial(-4122)

# New row 8: line7
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "line7"
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11
$ws.Range("E8").Value = $false

# New row 9: line8
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "line8"
$ws.Range("C9").Value = 16
$ws.Range("D9").Value = 9
$ws.Range("E9").Value = $true

# Update the shifted extr1..extr8 rows (now rows 10-17) with their new
# index (column A) and new C/D/E values per the target data.
$ws.Range("A10").Value = 8
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 12
$ws.Range("E10").Value = $true

$ws.Range("A11").Value = 9
$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 9
$ws.Range("E11").Value = $true

$ws.Range("A12").Value = 10
$ws.Range("C12").Value = 10
$ws.Range("D12").Value = 11
$ws.Range("E12").Value = $false

$ws.Range("A13").Value = 11
$ws.Range("C13").Value = 7
$ws.Range("D13").Value = 8
$ws.Range("E13").Value = $false

$ws.Range("A14").Value = 12
$ws.Range("C14").Value = 9
$ws.Range("D14").Value = 11
$ws.Range("E14").Value = $true

$ws.Range("A15").Value = 13
$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 11
$ws.Range("E15").Value = $true

$ws.Range("A16").Value = 14
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = $true

$ws.Range("A17").Value = 15
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = $false
